$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new data rows to the table with the new present list entries,
# and update the existing data row's values, matching the order the
# author actually typed them in (this affects shared-string ordering).
$ws.Range("A2").Value = "010"
$ws.Range("B2").Value = "בלה בלה"
$ws.Range("A3").Value = "050"
$ws.Range("A4").Value = "100"
$ws.Range("B4").Value = "בלה בלה בלה"
$ws.Range("B3").Value = "בלה"

# Expand the table range to include the new rows.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:B4"))

$ws.Range("B6").Select()
